$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 16:35"

# Refresh per-country COVID-19 stats and re-sort order (Casos totales desc).
# The underlying data source was refreshed and the table re-sorted by "Casos
# totales" (col B) descending, so both values and the country occupying a
# given rank can change.
# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 3773089
$ws.Range("C4").Value = 3077
$ws.Range("D4").Value = 1741626
$ws.Range("E4").Value = 1889358
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 41
$ws.Range("H4").Value = 142105

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1054247
$ws.Range("C6").Value = 13790
$ws.Range("D6").Value = 662652
$ws.Range("E6").Value = 365119
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 191
$ws.Range("H6").Value = 26476

# Row 13: Reino Unido
$ws.Range("A13").Value = "Reino Unido"
$ws.Range("B13").Value = 294066
$ws.Range("C13").Value = 827
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 45233

# Row 19: Alemania
$ws.Range("A19").Value = "Alemania"
$ws.Range("B19").Value = 202416
$ws.Range("C19").Value = 71
$ws.Range("D19").Value = 187500
$ws.Range("E19").Value = 5754
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 9162

# Row 41: Paises Bajos
$ws.Range("A41").Value = "Paises Bajos"
$ws.Range("B41").Value = 51581
$ws.Range("C41").Value = 127
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 6136

# Row 44: Portugal
$ws.Range("A44").Value = "Portugal"
$ws.Range("B44").Value = 48390
$ws.Range("C44").Value = 313
$ws.Range("D44").Value = 33153
$ws.Range("E44").Value = 13553
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 1684

# Row 56: Azerbaiyan
$ws.Range("A56").Value = "Azerbaiyan"
$ws.Range("B56").Value = 27133
$ws.Range("C56").Value = 497
$ws.Range("D56").Value = 18450
$ws.Range("E56").Value = 8334
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 8
$ws.Range("H56").Value = 349

# Row 79: Noruega
$ws.Range("A79").Value = "Noruega"
$ws.Range("B79").Value = 9028
$ws.Range("C79").Value = 3
$ws.Range("D79").Value = 8138
$ws.Range("E79").Value = 635
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 255

# Row 91: Tayikistan
$ws.Range("A91").Value = "Tayikistan"
$ws.Range("B91").Value = 6834
$ws.Range("C91").Value = 48
$ws.Range("D91").Value = 5529
$ws.Range("E91").Value = 1248
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 57

# Row 108: Zambia
$ws.Range("A108").Value = "Zambia"
$ws.Range("B108").Value = 2980
$ws.Range("C108").Value = 170
$ws.Range("D108").Value = 1462
$ws.Range("E108").Value = 1398
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 11
$ws.Range("H108").Value = 120

# Row 109: Maldivas
$ws.Range("A109").Value = "Maldivas"
$ws.Range("B109").Value = 2913
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 2340
$ws.Range("E109").Value = 558
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 15

# Row 110: Malaui
$ws.Range("A110").Value = "Malaui"
$ws.Range("B110").Value = 2810
$ws.Range("C110").Value = 5
$ws.Range("D110").Value = 1111
$ws.Range("E110").Value = 1644
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 55

# Row 112: Sri Lanka
$ws.Range("A112").Value = "Sri Lanka"
$ws.Range("B112").Value = 2703
$ws.Range("C112").Value = 6
$ws.Range("D112").Value = 2023
$ws.Range("E112").Value = 669
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 11

# Row 115: Mali
$ws.Range("A115").Value = "Mali"
$ws.Range("B115").Value = 2472
$ws.Range("C115").Value = 5
$ws.Range("D115").Value = 1809
$ws.Range("E115").Value = 542
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 121

# Row 134: Mozambique
$ws.Range("A134").Value = "Mozambique"
$ws.Range("B134").Value = 1435
$ws.Range("C134").Value = 33
$ws.Range("D134").Value = 408
$ws.Range("E134").Value = 1017
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 10

# Row 135: Zimbabue
$ws.Range("A135").Value = "Zimbabue"
$ws.Range("B135").Value = 1420
$ws.Range("C135").Value = 0
$ws.Range("D135").Value = 438
$ws.Range("E135").Value = 958
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 24

# Row 138: Namibia
$ws.Range("A138").Value = "Namibia"
$ws.Range("B138").Value = 1203
$ws.Range("C138").Value = 125
$ws.Range("D138").Value = 32
$ws.Range("E138").Value = 1169
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 2

# Row 139: Letonia
$ws.Range("A139").Value = "Letonia"
$ws.Range("B139").Value = 1189
$ws.Range("C139").Value = 4
$ws.Range("D139").Value = 1022
$ws.Range("E139").Value = 136
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 31

# Row 140: Niger
$ws.Range("A140").Value = "Niger"
$ws.Range("B140").Value = 1102
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 1013
$ws.Range("E140").Value = 20
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 69

# Row 141: Liberia
$ws.Range("A141").Value = "Liberia"
$ws.Range("B141").Value = 1088
$ws.Range("C141").Value = 3
$ws.Range("D141").Value = 519
$ws.Range("E141").Value = 499
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 70
